# Apply cryptos list updates (prices, 1h volume %, and a 3-row reorder
# around PancakeSwap / PEPE / Fetch.AI) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.184.33"
$ws.Range("E2").Value = "  +3.46%  "
$ws.Range("D3").Value = "2.313.13"
$ws.Range("E3").Value = "  +2.11%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.71"
$ws.Range("E5").Value = "  +4.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.54"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").Value = "2.333.40"
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("E10").Value = "  +8.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("E12").Value = "  +7.80%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.99"
$ws.Range("E14").Value = "  +4.32%  "
$ws.Range("D15").Value = "2.724.02"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").Value = "56.334.72"
$ws.Range("E16").Value = "  +3.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("E17").Value = "  +4.75%  "
$ws.Range("D18").Value = "2.317.16"
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("E19").Value = "  +2.92%  "
$ws.Range("E20").Value = "  +3.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.26"
$ws.Range("E21").Value = "  +6.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.63"
$ws.Range("E22").Value = "  +4.79%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.70"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.993"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("E26").Value = "  +6.22%  "
$ws.Range("E27").Value = "  +4.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.60"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.20"
$ws.Range("E29").Value = "  +10.79%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.70"
$ws.Range("E30").Value = "  +5.79%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0730"
$ws.Range("E31").Value = "  +5.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.26"
$ws.Range("E32").Value = "  +4.80%  "
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("E36").Value = "  +5.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.923"
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("E38").Value = "  +8.18%  "
$ws.Range("E39").Value = "  +8.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.48"
$ws.Range("E40").Value = "  +4.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.382"
$ws.Range("E41").Value = "  +1.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.58"
$ws.Range("E42").Value = "  +12.41%  "
$ws.Range("E43").Value = "  +6.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "278.64"
$ws.Range("E44").Value = "  +15.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.08"
$ws.Range("E45").Value = "  +5.70%  "
$ws.Range("E46").Value = "  +3.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0928"
$ws.Range("E47").Value = "  +3.72%  "
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("E49").Value = "  +2.10%  "
$ws.Range("E50").Value = "  +5.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.96"
$ws.Range("E51").Value = "  +5.15%  "
